$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Points Completed" value for Sprint 2
$ws.Range("C3").Value = 28

# Add the new Sprint 3 row
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 34

# Update the active selection as recorded for this edit
$ws.Range("C10").Select()

$wb.Save()
